$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "27.658.62"
Set-TextValue "E2" "  -2.51%  "
Set-TextValue "D3" "1.847.15"
Set-TextValue "E3" "  -0.85%  "
Set-TextValue "E4" "  -0.25%  "
Set-TextValue "D5" "314.43"
Set-TextValue "E5" "  -1.65%  "
Set-TextValue "E6" "  -0.16%  "
Set-TextValue "D7" "0.4268"
Set-TextValue "D8" "0.3674"
Set-TextValue "E8" "  -1.45%  "
Set-TextValue "D9" "44.83"
Set-TextValue "E9" "  -0.68%  "
Set-TextValue "D10" "0.07249"
Set-TextValue "E10" "  -3.94%  "
Set-TextValue "D11" "0.9011"
Set-TextValue "E11" "  -3.96%  "
Set-TextValue "D12" "20.74"
Set-TextValue "E12" "  -2.85%  "
Set-TextValue "D13" "1.876.22"
Set-TextValue "E13" "  -0.90%  "
Set-TextValue "D14" "6.580"
Set-TextValue "E14" "  -2.03%  "
Set-TextValue "D15" "5.356"
Set-TextValue "E15" "  -2.17%  "
Set-TextValue "D16" "0.06843"
Set-TextValue "E16" "  -0.59%  "
Set-TextValue "D17" "1.002"
Set-TextValue "E17" "  -0.30%  "
Set-TextValue "D18" "77.73"
Set-TextValue "E18" "  -5.49%  "
Set-TextValue "E20" "  -0.29%  "
Set-TextValue "E21" "  -3.49%  "
Set-TextValue "D22" "27.638.37"
Set-TextValue "E22" "  -2.55%  "
Set-TextValue "D23" "4.961"
Set-TextValue "E23" "  -3.18%  "
Set-TextValue "E24" "  -0.64%  "
Set-TextValue "D25" "2.093.01"
Set-TextValue "E25" "  -1.14%  "
Set-TextValue "D26" "2.046"
Set-TextValue "E26" "  +0.23%  "
Set-TextValue "D27" "154.05"
Set-TextValue "E27" "  -0.53%  "
Set-TextValue "D28" "18.23"
Set-TextValue "E28" "  -1.19%  "
Set-TextValue "D29" "5.251"
Set-TextValue "D30" "1.836"
Set-TextValue "E30" "  +5.74%  "
Set-TextValue "E31" "  -3.34%  "
Set-TextValue "D32" "0.08898"
Set-TextValue "E32" "  -2.27%  "
Set-TextValue "D33" "0.7706"
Set-TextValue "E33" "  -4.37%  "
Set-TextValue "D34" "4.561"
Set-TextValue "E34" "  -6.66%  "
Set-TextValue "D35" "2.917"
Set-TextValue "E35" "  -0.98%  "
Set-TextValue "D36" "1.085"
Set-TextValue "E36" "  -7.94%  "
Set-TextValue "E37" "  -0.27%  "
Set-TextValue "D38" "0.05384"
Set-TextValue "E38" "  -1.71%  "
Set-TextValue "D39" "1.095"
Set-TextValue "E39" "  -3.08%  "
Set-TextValue "B40" "VeChain"
Set-TextValue "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.01924"
Set-TextValue "E40" "  -1.59%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.948"
Set-TextValue "E41" "  -2.36%  "
Set-TextValue "D42" "0.5067"
Set-TextValue "E42" "  -3.79%  "
Set-TextValue "E43" "  -2.35%  "
Set-TextValue "D44" "6.763"
Set-TextValue "E44" "  -5.85%  "
Set-TextValue "D45" "8.271"
Set-TextValue "E45" "  -6.89%  "
Set-TextValue "D46" "0.06633"
Set-TextValue "E46" "  -2.16%  "
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.4724"
Set-TextValue "E47" "  -3.79%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "10.35"
Set-TextValue "E48" "  -2.24%  "
Set-TextValue "D49" "105.23"
Set-TextValue "E49" "  -2.51%  "
Set-TextValue "E50" "  -0.20%  "
Set-TextValue "E51" "  -2.96%  "
